$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test 4")
$ws.Activate()

$ws.Range("C7").Formula = "=24/B5"

$ws.Range("C8").Select()
